$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Cypher query text in cell B2 (WebExcel / CDS "Participants" query)
$newQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['BW']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

$ws.Range("B2").Value = $newQuery

# The longer query text now wraps across more lines - grow row 2 to fit
$ws.Rows.Item(2).RowHeight = 279

# Move the active selection from B3 down to B4 and scroll the view back to the top
$ws.Range("B4").Select() | Out-Null
